$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2932.1538
$ws.Range("I62").Value = 2979.6667
$ws.Range("J62").Value = 2825.25
$ws.Range("K62").Value = 2979.6667
$ws.Range("L62").Value = 2825.25
$ws.Range("M62").Value = -2355.6667
$ws.Range("N62").Value = -4073.25
$ws.Range("H65").Value = 2932.1538
$ws.Range("I65").Value = 2979.6667
$ws.Range("J65").Value = 2825.25
$ws.Range("K65").Value = 14898.3335
$ws.Range("L65").Value = 14126.25
$ws.Range("M65").Value = -11778.3335
$ws.Range("N65").Value = -20366.25
$ws.Range("H113").Value = 2942.4285
$ws.Range("I113").Value = 2899.25
$ws.Range("K113").Value = 2899.25
$ws.Range("M113").Value = 354.75
$ws.Range("H137").Value = 1858.5
$ws.Range("I137").Value = 1674.1111
$ws.Range("J137").Value = 2411.6667
$ws.Range("K137").Value = 5022.3333
$ws.Range("L137").Value = 7235.000100000001
$ws.Range("M137").Value = -2472.3333
$ws.Range("N137").Value = -12335.0001
$ws.Range("H138").Value = 2386602.5
$ws.Range("I138").Value = 5267140.5
$ws.Range("J138").Value = 7027.674
$ws.Range("K138").Value = 15801421.5
$ws.Range("L138").Value = 21083.022
$ws.Range("M138").Value = -15796281.5
$ws.Range("N138").Value = -31363.022
$ws.Range("H141").Value = 8247.916999999999
$ws.Range("I141").Value = 4258.696
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 12776.088
$ws.Range("L141").Value = 300000
$ws.Range("M141").Value = -7596.088
$ws.Range("N141").Value = -310360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 92482.73
$ws.Range("I102").Value = 1563.75
$ws.Range("K102").Value = 1563.75
$ws.Range("M102").Value = 58.25
$ws.Range("H132").Value = 2171.652
$ws.Range("I132").Value = 1830.6757
$ws.Range("J132").Value = 3573.4443
$ws.Range("K132").Value = 5492.0271
$ws.Range("L132").Value = 10720.3329
$ws.Range("M132").Value = -2962.0271
$ws.Range("N132").Value = -15780.3329
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1664.2941
$ws.Range("I99").Value = 1129.3
$ws.Range("J99").Value = 2428.5715
$ws.Range("K99").Value = 1129.3
$ws.Range("L99").Value = 2428.5715
$ws.Range("M99").Value = 368.7
$ws.Range("N99").Value = -5424.5715
$ws.Range("H105").Value = 2415.1924
$ws.Range("I105").Value = 2431.8
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2431.8
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -684.8000000000002
$ws.Range("N105").Value = -5494
$ws.Range("H134").Value = 3033.5178
$ws.Range("I134").Value = 1836.1708
$ws.Range("K134").Value = 5508.512400000001
$ws.Range("M134").Value = -2973.512400000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4169.5713
$ws.Range("I31").Value = 2898.8096
$ws.Range("J31").Value = 6075.7144
$ws.Range("K31").Value = 2898.8096
$ws.Range("L31").Value = 6075.7144
$ws.Range("M31").Value = -2603.8096
$ws.Range("N31").Value = -6665.7144
$ws.Range("H34").Value = 4169.5713
$ws.Range("I34").Value = 2898.8096
$ws.Range("J34").Value = 6075.7144
$ws.Range("K34").Value = 2898.8096
$ws.Range("L34").Value = 6075.7144
$ws.Range("M34").Value = -2696.8096
$ws.Range("N34").Value = -6479.7144
$ws.Range("H81").Value = 37000
$ws.Range("J81").Value = 37000
$ws.Range("L81").Value = 37000
$ws.Range("N81").Value = -38996
$ws.Range("H84").Value = 37000
$ws.Range("J84").Value = 37000
$ws.Range("L84").Value = 111000
$ws.Range("N84").Value = -120984
$ws.Range("H132").Value = 1424.4255
$ws.Range("I132").Value = 1368.6046
$ws.Range("K132").Value = 4105.8138
$ws.Range("M132").Value = -1575.8138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1287.4286
$ws.Range("I121").Value = 1020
$ws.Range("J121").Value = 1644
$ws.Range("K121").Value = 3060
$ws.Range("L121").Value = 4932
$ws.Range("M121").Value = -1750
$ws.Range("N121").Value = -7552
$ws.Range("H123").Value = 8400
$ws.Range("H131").Value = 34488144
$ws.Range("I131").Value = 25575
$ws.Range("J131").Value = 40002150
$ws.Range("K131").Value = 76725
$ws.Range("L131").Value = 120006450
$ws.Range("M131").Value = -71685
$ws.Range("N131").Value = -120016530
$ws.Range("H133").Value = 1970
$ws.Range("J133").Value = 3800
$ws.Range("L133").Value = 11400
$ws.Range("N133").Value = -21520
$ws.Range("H134").Value = 3890.4092
$ws.Range("I134").Value = 2536.8125
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 7610.4375
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = -2540.4375
$ws.Range("N134").Value = -32640
$ws.Range("H137").Value = 20838970
$ws.Range("I137").Value = 2676.6667
$ws.Range("J137").Value = 25647344
$ws.Range("K137").Value = 8030.000100000001
$ws.Range("L137").Value = 76942032
$ws.Range("M137").Value = -2930.000100000001
$ws.Range("N137").Value = -76952232
$ws.Range("H139").Value = 1804.3846
$ws.Range("I139").Value = 1442.4
$ws.Range("J139").Value = 3011
$ws.Range("K139").Value = 4327.200000000001
$ws.Range("L139").Value = 9033
$ws.Range("M139").Value = 812.7999999999993
$ws.Range("N139").Value = -19313
$ws.Range("H140").Value = 2834.818
$ws.Range("I140").Value = 2315
$ws.Range("K140").Value = 6945
$ws.Range("M140").Value = -1765
$ws.Range("H141").Value = 4927.6665
$ws.Range("I141").Value = 4927.6665
$ws.Range("K141").Value = 14782.9995
$ws.Range("M141").Value = -9602.999500000002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 29757
$ws.Range("J64").Value = 29757
$ws.Range("L64").Value = 29757
$ws.Range("N64").Value = -30253
$ws.Range("H67").Value = 29757
$ws.Range("J67").Value = 29757
$ws.Range("L67").Value = 29757
$ws.Range("N67").Value = -31473
$ws.Range("H132").Value = 1603.275
$ws.Range("I132").Value = 1338.4482
$ws.Range("K132").Value = 4015.3446
$ws.Range("M132").Value = -1485.3446
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4545.1577
$ws.Range("I7").Value = 4504.077
$ws.Range("J7").Value = 4634.1665
$ws.Range("K7").Value = 4504.077
$ws.Range("L7").Value = 4634.1665
$ws.Range("M7").Value = -4392.077
$ws.Range("N7").Value = -4858.1665
$ws.Range("H126").Value = 4545.1577
$ws.Range("I126").Value = 4504.077
$ws.Range("J126").Value = 4634.1665
$ws.Range("K126").Value = 13512.231
$ws.Range("L126").Value = 13902.4995
$ws.Range("M126").Value = -11042.231
$ws.Range("N126").Value = -18842.4995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H136").Value = 1824.2
$ws.Range("I136").Value = 1860.6333
$ws.Range("J136").Value = 1714.9
$ws.Range("K136").Value = 5581.8999
$ws.Range("L136").Value = 5144.700000000001
$ws.Range("M136").Value = -3031.8999
$ws.Range("N136").Value = -10244.7
